# Prefix each "Step..." / label command name in column A (rows 2..N) of the
# protocol worksheets with the worksheet's own name, so that command names
# are unique across the workbook (e.g. "Step4 Takeaway" -> "discount2 Step4 Takeaway").
#
# This applies to every worksheet whose name matches one of the known
# protocol/step sheets (price1, price2, discount1, discount2, free1, free2,
# nomoney1, nomoney2, noppv1, noppv2, card1, card2, nosex1, nosex2,
# offtopic1, offtopic2, real1, real2, voice1, voice2, customyes1,
# customyes2, customno1, customno2, done1, done2, cumcontrol, dickpic,
# boosters). The first six worksheets (Jack HollywoodJourney, MeetupRedirect,
# NRWaves, PersonalJack Hollywood, PositiveSpin, ReEngagement) are left
# untouched.

$wb = $excel.ActiveWorkbook

$targetSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    if ($ws -eq $null) { continue }

    $prefix = $sheetName + " "

    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2

        if ($current -eq $null) { continue }
        if ($current -eq "") { continue }

        # Avoid double-prefixing if already prefixed.
        if ($current.StartsWith($prefix)) { continue }

        $cell.Value = $prefix + $current
    }
}
